# Auto-generated script to update cryptos.xlsx Price and Volume(1h) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '64.281.08'
$ws.Cells.Item(2, 5).Value = '  +0.33%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.496.99'
$ws.Cells.Item(3, 5).Value = '  -0.56%  '
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '588.30'
$ws.Cells.Item(5, 5).Value = '  +0.35%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '134.28'
$ws.Cells.Item(6, 5).Value = '  +0.77%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 5).Value = '  +0.26%  '
$ws.Cells.Item(9, 5).Value = '  +0.28%  '
$ws.Cells.Item(10, 5).Value = '  +2.43%  '
$ws.Cells.Item(11, 5).Value = '  +2.47%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '4.092.48'
$ws.Cells.Item(12, 5).Value = '  -0.35%  '
$ws.Cells.Item(13, 5).Value = '  +1.37%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.0000180'
$ws.Cells.Item(14, 5).Value = '  +1.27%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '3.497.96'
$ws.Cells.Item(15, 5).Value = '  -0.35%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '64.364.06'
$ws.Cells.Item(16, 5).Value = '  +0.44%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '25.62'
$ws.Cells.Item(17, 5).Value = '  -6.71%  '
$ws.Cells.Item(18, 5).Value = '  +0.46%  '
$ws.Cells.Item(19, 5).Value = '  +2.43%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '13.53'
$ws.Cells.Item(20, 5).Value = '  -2.70%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '393.12'
$ws.Cells.Item(21, 5).Value = '  +2.57%  '
$ws.Cells.Item(22, 5).Value = '  +0.01%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.636.80'
$ws.Cells.Item(23, 5).Value = '  -0.48%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '74.58'
$ws.Cells.Item(24, 5).Value = '  +0.93%  '
$ws.Cells.Item(25, 5).Value = '  +0.10%  '
$ws.Cells.Item(26, 5).Value = '  +1.33%  '
$ws.Cells.Item(27, 5).Value = '  +1.11%  '
$ws.Cells.Item(28, 5).Value = '  +0.06%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '7.35'
$ws.Cells.Item(29, 5).Value = '  -1.73%  '
$ws.Cells.Item(30, 5).Value = '  +0.36%  '
$ws.Cells.Item(31, 5).Value = '  -1.62%  '
$ws.Cells.Item(32, 5).Value = '  -6.28%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.519.88'
$ws.Cells.Item(33, 5).Value = '  -0.11%  '
$ws.Cells.Item(34, 5).Value = '  +5.27%  '
$ws.Cells.Item(35, 5).Value = '  +0.03%  '
$ws.Cells.Item(36, 5).Value = '  -0.37%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '5.14'
$ws.Cells.Item(37, 5).Value = '  -4.13%  '
$ws.Cells.Item(38, 5).Value = '  -0.68%  '
$ws.Cells.Item(39, 5).Value = '  -1.22%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '167.07'
$ws.Cells.Item(40, 5).Value = '  +3.97%  '
$ws.Cells.Item(41, 5).Value = '  -0.48%  '
$ws.Cells.Item(42, 5).Value = '  -0.46%  '
$ws.Cells.Item(43, 5).Value = '  +0.07%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '25.42'
$ws.Cells.Item(44, 5).Value = '  -4.41%  '
$ws.Cells.Item(45, 5).Value = '  -0.29%  '
$ws.Cells.Item(46, 5).Value = '  +2.63%  '
$ws.Cells.Item(47, 5).Value = '  -3.95%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '2.465.75'
$ws.Cells.Item(48, 5).Value = '  -0.35%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '6.75'
$ws.Cells.Item(49, 5).Value = '  -0.74%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.891'
$ws.Cells.Item(50, 5).Value = '  -1.87%  '
$ws.Cells.Item(51, 5).Value = '  -1.27%  '
